# Apply the "New merchantpage and Test page added" edit to the
# RunManager sheet (3rd worksheet) of the CloudPos workbook.
#
# Changes:
#  - Rows 2-5, column C ("Execute") switch from "yes" to "NO"
#  - A brand new row 6 is appended for a new test case
#    "validateMerchantPageOnSystemLogin" with Execute = "yes",
#    following the same layout/format as the existing rows.
#  - The active selection moves to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RunManager")

# Add the new test row first, mirroring the formatting of the row above
# it (single-call Copy keeps the destination's formats/values in sync
# without disturbing any other existing cells).
$ws.Range("A5:D5").Copy($ws.Range("A6:D6"))
$excel.CutCopyMode = 0

$ws.Range("A6").Formula = "validateMerchantPageOnSystemLogin"
$ws.Range("B6").Formula = "To Check Whether"
# Keep the leading apostrophe so the text-quote-prefix formatting of the
# copied cell (column is formatted as Text) is preserved.
$ws.Range("C6").Formula = "'yes"
$ws.Range("D6").Formula = "'1"

# Disable the previously active test cases. The leading apostrophe keeps
# these text cells' existing quote-prefixed "Text" formatting intact.
$ws.Range("C2").Formula = "'NO"
$ws.Range("C3").Formula = "'NO"
$ws.Range("C4").Formula = "'NO"
$ws.Range("C5").Formula = "'NO"

# Match the saved selection state recorded in the workbook.
[void]$ws.Range("D5").Select()
